$wb = $excel.ActiveWorkbook

# --- Rename sheets: shift "Activité N" numbering by +4 (5->9, 6->10, 7->11) ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$ws1.Name = "Activité 9"
$ws2.Name = "Activité 10"
$ws3.Name = "Activité 11"

# --- Move the active tab / selected sheet from the 1st sheet to the 3rd sheet ---
$ws3.Activate()

# --- Fix the header/footer font-style label on every sheet: "Regular" -> "Normal" ---
foreach ($ws in @($ws1, $ws2, $ws3)) {
    $ps = $ws.PageSetup
    $ps.CenterHeader = '&"Times New Roman,Normal"&12&A'
    $ps.CenterFooter = '&"Times New Roman,Normal"&12Page &P'
}

# --- Default column width tweaks (best effort) ---
$ws1.StandardWidth = 11.55078125
$ws2.StandardWidth = 11.5703125
$ws3.StandardWidth = 11.5703125
